# Edit: split " (romák, melegek, migránsok" / ")" runs so that a new
# "bevándorlók/" run is inserted right before the "_GoBack" bookmark, turning
#   " (romák, melegek, migránsok" + ")"
# into
#   " (romák, melegek, " + "bevándorlók/" + <bookmark> + "migránsok)"
#
# The Word object model in this runtime merges adjacent, identically
# formatted runs whenever text is inserted via Range.InsertBefore /
# InsertAfter, so the only reliable way to create a genuinely separate
# sibling <w:r> (as required by the target OOXML) is to round-trip the
# document part through Range.XML()/Range.InsertXML().

$d = $word.ActiveDocument
$full = $d.Content

# 1) Pull the whole package XML (every part) out of the document.
$pkgXml = $full.XML($true)

# 2) Isolate just the /word/document.xml part's contents.
$partPattern = '<pkg:part pkg:name="/word/document.xml"[^>]*><pkg:xmlData>(.*?)</pkg:xmlData></pkg:part>'
$partMatch = [regex]::Match($pkgXml, $partPattern, [System.Text.RegularExpressions.RegexOptions]::Singleline)
if (-not $partMatch.Success) {
    throw "Could not locate /word/document.xml part inside package XML"
}
$docXml = $partMatch.Groups[1].Value

# 3) Apply the precise textual edit described by the diff.
$oldFragment = '<w:r><w:rPr><w:rFonts w:ascii="Calisto MT" w:hAnsi="Calisto MT"/></w:rPr><w:t xml:space="preserve"> (romák, melegek, migránsok</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:rFonts w:ascii="Calisto MT" w:hAnsi="Calisto MT"/></w:rPr><w:t>)</w:t></w:r>'
$newFragment = '<w:r><w:rPr><w:rFonts w:ascii="Calisto MT" w:hAnsi="Calisto MT"/></w:rPr><w:t xml:space="preserve"> (romák, melegek, </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calisto MT" w:hAnsi="Calisto MT"/></w:rPr><w:t>bevándorlók/</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:rFonts w:ascii="Calisto MT" w:hAnsi="Calisto MT"/></w:rPr><w:t>migránsok)</w:t></w:r>'

if (-not $docXml.Contains($oldFragment)) {
    throw "Expected original run sequence was not found in document.xml"
}
$newDocXml = $docXml.Replace($oldFragment, $newFragment)

# 4) Re-wrap just the document.xml part in a minimal package and push it
#    back into the document. Limiting the package to a single part avoids
#    needlessly touching the other parts (styles/media/etc.) on re-import.
$newPkgXml = '<?xml version="1.0" standalone="yes"?>' + "`n" +
    '<?mso-application progid="Word.Document"?>' + "`n" +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' + $newDocXml + '</pkg:xmlData></pkg:part></pkg:package>'

$full.InsertXML($newPkgXml)

Write-Host "Edit applied."
